$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")

# Row 5
$ws.Range("H5").Value = 125000420
$ws.Range("I5").Value = 123.5
$ws.Range("K5").Value = 123.5
$ws.Range("M5").Value = -8.5

# Row 17
$ws.Range("H17").Value = 1777.3334
$ws.Range("J17").Value = 1923.5555
$ws.Range("L17").Value = 5770.666499999999
$ws.Range("N17").Value = -6106.666499999999

# Row 28
$ws.Range("H28").Value = 3983
$ws.Range("I28").Value = 3198.4348
$ws.Range("K28").Value = 3198.4348
$ws.Range("M28").Value = -2713.4348

# Row 33
$ws.Range("H33").Value = 210.82353
$ws.Range("I33").Value = 155.3077
$ws.Range("K33").Value = 155.3077
$ws.Range("M33").Value = 73.69229999999999

# Row 62
$ws.Range("H62").Value = 10419634
$ws.Range("I62").Value = 13891567
$ws.Range("J62").Value = 3835.1667
$ws.Range("K62").Value = 13891567
$ws.Range("L62").Value = 3835.1667
$ws.Range("M62").Value = -13890943
$ws.Range("N62").Value = -5083.1667

# Row 65
$ws.Range("H65").Value = 10419634
$ws.Range("I65").Value = 13891567
$ws.Range("J65").Value = 3835.1667
$ws.Range("K65").Value = 69457835
$ws.Range("L65").Value = 19175.8335
$ws.Range("M65").Value = -69454715
$ws.Range("N65").Value = -25415.8335

# Row 106
$ws.Range("H106").Value = 4789682.5
$ws.Range("I106").Value = 5687075
$ws.Range("J106").Value = 3589.6667
$ws.Range("K106").Value = 5687075
$ws.Range("L106").Value = 3589.6667
$ws.Range("M106").Value = -5686444
$ws.Range("N106").Value = -4851.6667

# Row 125
$ws.Range("H125").Value = 4129.5
$ws.Range("I125").Value = 2367.3333
$ws.Range("J125").Value = 4884.7144
$ws.Range("K125").Value = 21305.9997
$ws.Range("L125").Value = 43962.4296
$ws.Range("M125").Value = -18845.9997
$ws.Range("N125").Value = -48882.4296

# Row 127
$ws.Range("H127").Value = 620.8333
$ws.Range("I127").Value = 518.2727
$ws.Range("J127").Value = 1749
$ws.Range("K127").Value = 1554.8181
$ws.Range("L127").Value = 5247
$ws.Range("M127").Value = 3405.1819
$ws.Range("N127").Value = -15167

# Row 137
$ws.Range("H137").Value = 1004.75555
$ws.Range("I137").Value = 814.4857
$ws.Range("J137").Value = 1670.7
$ws.Range("K137").Value = 2443.4571
$ws.Range("L137").Value = 5012.1
$ws.Range("M137").Value = 106.5429000000004
$ws.Range("N137").Value = -10112.1

# Row 141
$ws.Range("H141").Value = 6938.5454
$ws.Range("I141").Value = 7247.5
$ws.Range("K141").Value = 21742.5
$ws.Range("M141").Value = -16562.5

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 1020.96155
$ws.Range("J2").Value = 1779.8
$ws.Range("L2").Value = 1779.8
$ws.Range("N2").Value = -2005.8

# Row 63
$ws.Range("H63").Value = 766.1111
$ws.Range("I63").Value = 415.14285
$ws.Range("K63").Value = 415.14285
$ws.Range("M63").Value = 270.85715

# Row 66
$ws.Range("H66").Value = 766.1111
$ws.Range("I66").Value = 415.14285
$ws.Range("K66").Value = 2075.71425
$ws.Range("M66").Value = 1356.28575

# Row 97
$ws.Range("H97").Value = 2629.353
$ws.Range("I97").Value = 562.72
$ws.Range("K97").Value = 562.72
$ws.Range("M97").Value = -66.72000000000003

# Row 116
$ws.Range("H116").Value = 1020.96155
$ws.Range("J116").Value = 1779.8
$ws.Range("L116").Value = 1779.8
$ws.Range("N116").Value = -6367.8

# Row 132
$ws.Range("H132").Value = 3794.7307
$ws.Range("I132").Value = 3855.1
$ws.Range("J132").Value = 3593.5
$ws.Range("K132").Value = 11565.3
$ws.Range("L132").Value = 10780.5
$ws.Range("M132").Value = -9035.299999999999
$ws.Range("N132").Value = -15840.5

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 1020.96155
$ws.Range("J3").Value = 1779.8
$ws.Range("L3").Value = 1779.8
$ws.Range("N3").Value = -2007.8

# Row 134
$ws.Range("H134").Value = 7187.8823
$ws.Range("I134").Value = 7168.5947
$ws.Range("J134").Value = 7238.857
$ws.Range("K134").Value = 21505.7841
$ws.Range("L134").Value = 21716.571
$ws.Range("M134").Value = -18970.7841
$ws.Range("N134").Value = -26786.571

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")

# Row 16
$ws.Range("H16").Value = 1651.9166
$ws.Range("I16").Value = 1347.5454
$ws.Range("K16").Value = 1347.5454
$ws.Range("M16").Value = -1060.5454

# Row 22
$ws.Range("H22").Value = 438.8421
$ws.Range("I22").Value = 327.72726
$ws.Range("J22").Value = 591.625
$ws.Range("K22").Value = 327.72726
$ws.Range("L22").Value = 591.625
$ws.Range("M22").Value = 22.27274
$ws.Range("N22").Value = -1291.625

# Row 31
$ws.Range("H31").Value = 13114.143
$ws.Range("I31").Value = 16608.96
$ws.Range("J31").Value = 3018
$ws.Range("K31").Value = 16608.96
$ws.Range("L31").Value = 3018
$ws.Range("M31").Value = -16313.96
$ws.Range("N31").Value = -3608

# Row 34
$ws.Range("H34").Value = 13114.143
$ws.Range("I34").Value = 16608.96
$ws.Range("J34").Value = 3018
$ws.Range("K34").Value = 16608.96
$ws.Range("L34").Value = 3018
$ws.Range("M34").Value = -16406.96
$ws.Range("N34").Value = -3422

# Row 58
$ws.Range("H58").Value = 2536.7273
$ws.Range("J58").Value = 3215.4
$ws.Range("L58").Value = 3215.4
$ws.Range("N58").Value = -3621.4

# Row 113
$ws.Range("H113").Value = 1651.9166
$ws.Range("I113").Value = 1347.5454
$ws.Range("K113").Value = 1347.5454
$ws.Range("M113").Value = 822.4546

# Row 132
$ws.Range("H132").Value = 2470.0476
$ws.Range("I132").Value = 1800.4193
$ws.Range("J132").Value = 4357.1816
$ws.Range("K132").Value = 5401.257900000001
$ws.Range("L132").Value = 13071.5448
$ws.Range("M132").Value = -2871.257900000001
$ws.Range("N132").Value = -18131.5448

# Row 136
$ws.Range("H136").Value = 2536.7273
$ws.Range("J136").Value = 3215.4
$ws.Range("L136").Value = 9646.200000000001
$ws.Range("N136").Value = -14746.2

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")

# Row 4
$ws.Range("H4").Value = 234656460
$ws.Range("I4").Value = 158701170
$ws.Range("K4").Value = 476103510
$ws.Range("M4").Value = -476103398

# Row 5
$ws.Range("H5").Value = 492.7037
$ws.Range("I5").Value = 447.85715
$ws.Range("J5").Value = 649.6667
$ws.Range("K5").Value = 1343.57145
$ws.Range("L5").Value = 1949.0001
$ws.Range("M5").Value = -1231.57145
$ws.Range("N5").Value = -2173.0001

# Row 40
$ws.Range("H40").Value = 82.75
$ws.Range("I40").Value = 70.333336
$ws.Range("K40").Value = 281.333344
$ws.Range("M40").Value = -212.333344

# Row 61
$ws.Range("H61").Value = 103.71429
$ws.Range("I61").Value = 104.333336
$ws.Range("K61").Value = 313.000008
$ws.Range("M61").Value = -98.00000799999998

# Row 132
$ws.Range("H132").Value = 3995.28
$ws.Range("I132").Value = 4372.5
$ws.Range("J132").Value = 3962.4783
$ws.Range("K132").Value = 39352.5
$ws.Range("L132").Value = 35662.3047
$ws.Range("M132").Value = -36822.5
$ws.Range("N132").Value = -40722.3047

# Row 135
$ws.Range("H135").Value = 492.7037
$ws.Range("I135").Value = 447.85715
$ws.Range("J135").Value = 649.6667
$ws.Range("K135").Value = 4030.71435
$ws.Range("L135").Value = 5847.0003
$ws.Range("M135").Value = -1495.71435
$ws.Range("N135").Value = -10917.0003

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")

# Row 132
$ws.Range("H132").Value = 2338.806
$ws.Range("I132").Value = 1805.1132
$ws.Range("K132").Value = 5415.3396
$ws.Range("M132").Value = -2885.3396

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")

# Row 40
$ws.Range("H40").Value = 2403.3044
$ws.Range("I40").Value = 2206.125
$ws.Range("K40").Value = 2206.125
$ws.Range("M40").Value = -2070.125

# Row 68
$ws.Range("H68").Value = 2279.8235
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Row 71
$ws.Range("H71").Value = 2279.8235
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Row 74
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# Row 77
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")

# Row 31
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
